$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 611
$ws1.Range("F4").Value = 538
$ws1.Range("F5").Value = 522
$ws1.Range("F6").Value = 291
$ws1.Range("F7").Value = 2631
$ws1.Range("F9").Value = 7230
$ws1.Range("F11").Value = 451
$ws1.Range("F13").Value = 168

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 611
$ws4.Range("F4").Value = 538
$ws4.Range("F5").Value = 522
$ws4.Range("F6").Value = 291
$ws4.Range("F9").Value = 2631
$ws4.Range("F11").Value = 7230
$ws4.Range("F13").Value = 451
$ws4.Range("F17").Value = 168
